$wb = $excel.ActiveWorkbook

# --- MAIN_CONTROLLER sheet ---
$main = $wb.Worksheets.Item("MAIN_CONTROLLER")

# Row for "FOS" (Si_No=1) RunStatus changes from "N" to "n"
$main.Range("B2").Value = "n"

# Row for "FOS" (Si_No=4, the UW row) gets highlighted red (same "N" value, add fill)
$main.Range("B5").Interior.Color = 255

# Update the active selection to B2 on this sheet
$main.Range("B2").Select()

# --- DATASHEET sheet ---
$data = $wb.Worksheets.Item("DATASHEET")

# Ishine row ExplicityWait (F4) changes from 20 to 30
$data.Range("F4").Value = 30

# Update the active selection to F4 on this sheet (also keeps it the active/tabbed sheet)
$data.Range("F4").Select()
